$wb = $excel.ActiveWorkbook

# Delete the stray row 16 ("Sheet" / 3 / 4) on the optimization_parameters sheet.
$wsOpt = $wb.Worksheets.Item("optimization_parameters")
$wsOpt.Rows.Item(16).Delete()

